$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update the two measured input values (C5, C7). The dependent D-column
# formulas (=Cx/C5) recalc automatically from these, which also ripples
# into the chart that plots column D.
$ws.Range("C5").Value = 527.2084
$ws.Range("C7").Value = 1091.297

$excel.CalculateFullRebuild()

# Update the active selection on the sheet to C6
$ws.Range("C6").Select()
